# Apply updated "想去人数" (F column) counts across the workbook's sheets.
# This mirrors the commit "Update gh-pages to output generated at 456a3b4",
# which refreshed the scraped attendance numbers for several rows.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    "F2"  = 3067
    "F3"  = 479
    "F4"  = 58
    "F5"  = 47
    "F6"  = 2
    "F8"  = 6
    "F9"  = 1053
    "F10" = 14846
    "F11" = 180
    "F12" = 142
    "F13" = 2
    "F14" = 5919
    "F15" = 606
    "F16" = 88
    "F18" = 86
    "F20" = 19
    "F23" = 813
    "F24" = 2955
    "F25" = 98
    "F26" = 10733
    "F28" = 77
    "F29" = 113
    "F30" = 3753
}
foreach ($cell in $sheet1Updates.Keys) {
    $ws1.Range($cell).Value = $sheet1Updates[$cell]
}

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 14

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    "F3"  = 3067
    "F4"  = 479
    "F5"  = 58
    "F6"  = 47
    "F7"  = 2
    "F9"  = 6
    "F10" = 1053
    "F11" = 14846
    "F12" = 180
    "F13" = 142
    "F14" = 2
    "F15" = 5919
    "F16" = 606
    "F17" = 88
    "F19" = 86
    "F21" = 19
    "F24" = 813
    "F25" = 2955
    "F26" = 98
    "F27" = 14
    "F28" = 10733
    "F30" = 77
    "F31" = 113
    "F32" = 3753
}
foreach ($cell in $sheet4Updates.Keys) {
    $ws4.Range($cell).Value = $sheet4Updates[$cell]
}
